$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "MoCoClf" (sheet3): add a comment to existing row 5, then append
# three new experiment rows (6, 7, 8).
# New shared strings introduced here, in first-use order, become indices
# 31, 32, 33 ("removed softmax before crossentropy", "similar results",
# "8630046_19").
# ---------------------------------------------------------------------------
$wsClf = $wb.Worksheets.Item("MoCoClf")

# Row 5 gains a comment in column L (reuses existing shared string
# "very similar to above").
$wsClf.Cells.Item(5, 12).Value = "very similar to above"

# Row 6
$wsClf.Cells.Item(6, 1).Value = 8641912
$wsClf.Cells.Item(6, 2).Value = "replace MoCo.fc -> 100 -> 4"
$wsClf.Cells.Item(6, 3).Value = 0.003
$wsClf.Cells.Item(6, 4).Value = "SGD"
$wsClf.Cells.Item(6, 5).Value = 128
$wsClf.Cells.Item(6, 6).Value = 1
$wsClf.Cells.Item(6, 7).Value = 50
$wsClf.Cells.Item(6, 8).Value = "8630046_39"
$wsClf.Cells.Item(6, 12).Value = "removed softmax before crossentropy"

# Row 7
$wsClf.Cells.Item(7, 1).Value = 8642199
$wsClf.Cells.Item(7, 2).Value = "replace MoCo.fc -> 100 -> 4"
$wsClf.Cells.Item(7, 3).Value = 0.0003
$wsClf.Cells.Item(7, 4).Value = "SGD"
$wsClf.Cells.Item(7, 5).Value = 128
$wsClf.Cells.Item(7, 6).Value = 1
$wsClf.Cells.Item(7, 7).Value = 50
$wsClf.Cells.Item(7, 8).Value = "8630046_39"
$wsClf.Cells.Item(7, 12).Value = "similar results"

# Row 8
$wsClf.Cells.Item(8, 1).Value = 8642862
$wsClf.Cells.Item(8, 2).Value = "replace MoCo.fc -> 100 -> 4"
$wsClf.Cells.Item(8, 3).Value = 0.0003
$wsClf.Cells.Item(8, 4).Value = "SGD"
$wsClf.Cells.Item(8, 5).Value = 128
$wsClf.Cells.Item(8, 6).Value = 1
$wsClf.Cells.Item(8, 7).Value = 50
$wsClf.Cells.Item(8, 8).Value = "8630046_19"

# ---------------------------------------------------------------------------
# Sheet "MoCo" (sheet2): append three new experiment rows (6, 7, 8).
# New shared strings introduced here become indices 34, 35, 36
# ("LabelMoCo", "converge to trivial solution", "increase memsize to
# 12800").
# ---------------------------------------------------------------------------
$wsMoCo = $wb.Worksheets.Item("MoCo")

# Row 6
$wsMoCo.Cells.Item(6, 1).Value = 8643007
$wsMoCo.Cells.Item(6, 2).Value = "LabelMoCo"
$wsMoCo.Cells.Item(6, 3).Value = 0.03
$wsMoCo.Cells.Item(6, 4).Value = "SGD"
$wsMoCo.Cells.Item(6, 5).Value = 128
$wsMoCo.Cells.Item(6, 6).Value = 1
$wsMoCo.Cells.Item(6, 7).Value = 100
$wsMoCo.Cells.Item(6, 8).Value = "converge to trivial solution"

# Row 7
$wsMoCo.Cells.Item(7, 1).Value = 8643330
$wsMoCo.Cells.Item(7, 2).Value = "LabelMoCo"
$wsMoCo.Cells.Item(7, 3).Value = 0.03
$wsMoCo.Cells.Item(7, 4).Value = "SGD"
$wsMoCo.Cells.Item(7, 5).Value = 128
$wsMoCo.Cells.Item(7, 6).Value = 1
$wsMoCo.Cells.Item(7, 7).Value = 100

# Row 8
$wsMoCo.Cells.Item(8, 1).Value = 8704485
$wsMoCo.Cells.Item(8, 2).Value = "LabelMoCo"
$wsMoCo.Cells.Item(8, 3).Value = 0.003
$wsMoCo.Cells.Item(8, 4).Value = "SGD"
$wsMoCo.Cells.Item(8, 5).Value = 128
$wsMoCo.Cells.Item(8, 6).Value = 1
$wsMoCo.Cells.Item(8, 7).Value = 100
$wsMoCo.Cells.Item(8, 8).Value = "increase memsize to 12800"

# ---------------------------------------------------------------------------
# Window / view state: reflect the new selections and zoom levels recorded
# in each sheet, finishing on the MoCoClf tab so it remains the active tab.
# ---------------------------------------------------------------------------
[void]$wsMoCo.Select()
$excel.ActiveWindow.Zoom = 180
[void]$wsMoCo.Range("I7").Select()

[void]$wsClf.Select()
$excel.ActiveWindow.Zoom = 181
[void]$wsClf.Range("A8").Select()

Write-Host "applied edits"
